# Generate Report for Handback
# - Update the "Ready for handoff" status to "Handback transform failed"
#   (this shared string is used on the Overview sheet rows for the second
#   file, as well as the per-locale sheets' Status column).
# - Populate the "Error Detail" column (P) on the zh-cn and de-de sheets
#   for the second file row with the handback/handoff filename mismatch
#   message.
# - Widen the "Error Detail" column (P) on the zh-cn and de-de sheets so
#   the new, longer error text is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 (a14a1f5d-...) Status / Path And Name columns
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-locale sheets: row 3 Status column (column C)
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Per-locale sheets: row 3 Error Detail column (column P)
$zhCnError = "Handback file name: ebpgtpxe.ynd is different with handoff file name: a14a1f5d-fe73-4bcf-bb9b-8664d3aab832.13cbc0f7c6739cd97325fb3aaed8bbfe8420dfef.zh-cn."
$deDeError = "Handback file name: ebpgtpxe.ynd is different with handoff file name: a14a1f5d-fe73-4bcf-bb9b-8664d3aab832.13cbc0f7c6739cd97325fb3aaed8bbfe8420dfef.de-de."

$wsZhCn.Range("P3").Value = $zhCnError
$wsDeDe.Range("P3").Value = $deDeError

# Widen the Error Detail column (P / column 16) on both locale sheets.
# NOTE: this runtime's ColumnWidth<->stored-width conversion has a small
# offset (character-width padding), so assigning 40 directly persists as
# ~40.83 in the saved XML. 39.1 is the COM-side value that round-trips to
# exactly width="40" in the OOXML, matching the target column width.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1

"Handback report generated"
